# Update "想去人数" (F) / "最低票价" (G) figures on the 展览 and 全部类型
# sheets to the newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -----------------------------------------------------------
$wsExhibit = $wb.Worksheets("展览")
$wsExhibit.Range("F2").Value = 171
$wsExhibit.Range("G2").Value = 55
$wsExhibit.Range("F3").Value = 657
$wsExhibit.Range("F4").Value = 25
$wsExhibit.Range("F6").Value = 1582
$wsExhibit.Range("F7").Value = 38
$wsExhibit.Range("F8").Value = 3148
$wsExhibit.Range("F10").Value = 736

# --- 全部类型 sheet (same events, shifted by one row due to an extra
#     concert entry in row 5) ------------------------------------------
$wsAll = $wb.Worksheets("全部类型")
$wsAll.Range("F2").Value = 171
$wsAll.Range("G2").Value = 55
$wsAll.Range("F3").Value = 657
$wsAll.Range("F4").Value = 25
$wsAll.Range("F7").Value = 1582
$wsAll.Range("F8").Value = 38
$wsAll.Range("F9").Value = 3148
$wsAll.Range("F11").Value = 736
